$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6429
$ws.Range("E2").Value = 342
$ws.Range("F2").Value = 342
$ws.Range("G2").Value = 359
$ws.Range("H2").Value = 269
$ws.Range("I2").Value = 269
$ws.Range("K2").Value = 3452
$ws.Range("L2").Value = 1263
$ws.Range("M2").Value = 2189
$ws.Range("N2").Value = 2189
$ws.Range("P2").Value = 173
$ws.Range("Q2").Value = 248
$ws.Range("R2").Value = -62
$ws.Range("S2").Value = -128
$ws.Range("T2").Value = 117
$ws.Range("U2").Value = 131
$ws.Range("V2").Value = 30
$ws.Range("W2").Value = 5.32
$ws.Range("X2").Value = 4.19
$ws.Range("Y2").Value = 13.08
$ws.Range("Z2").Value = 7.95
$ws.Range("AA2").Value = 57.69
$ws.Range("AB2").Value = 1186.65
$ws.Range("AC2").Value = 1561
$ws.Range("AD2").Value = 24.29
$ws.Range("AE2").Value = 12999
$ws.Range("AF2").Value = 2.92
$ws.Range("AG2").Value = 80
$ws.Range("AH2").Value = 0.21
$ws.Range("AI2").Value = 5
$ws.Range("AJ2").Value = 17261650
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 6942
$ws.Range("E3").Value = 390
$ws.Range("F3").Value = 403
$ws.Range("G3").Value = 404
$ws.Range("H3").Value = 288
$ws.Range("I3").Value = 288
$ws.Range("K3").Value = 3792
$ws.Range("L3").Value = 1335
$ws.Range("M3").Value = 2458
$ws.Range("N3").Value = 2458
$ws.Range("P3").Value = 173
$ws.Range("Q3").Value = 148
$ws.Range("R3").Value = -97
$ws.Range("S3").Value = -8
$ws.Range("T3").Value = 102
$ws.Range("U3").Value = 46
$ws.Range("V3").Value = 35
$ws.Range("W3").Value = 5.62
$ws.Range("X3").Value = 4.15
$ws.Range("Y3").Value = 12.41
$ws.Range("Z3").Value = 7.96
$ws.Range("AA3").Value = 54.31
$ws.Range("AB3").Value = 1341.02
$ws.Range("AC3").Value = 1670
$ws.Range("AD3").Value = 20.96
$ws.Range("AE3").Value = 14593
$ws.Range("AF3").Value = 2.4
$ws.Range("AG3").Value = 80
$ws.Range("AH3").Value = 0.23
$ws.Range("AI3").Value = 4.67
$ws.Range("AJ3").Value = 17261650
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 7356
$ws.Range("E4").Value = 422
$ws.Range("F4").Value = 422
$ws.Range("G4").Value = 449
$ws.Range("H4").Value = 332
$ws.Range("I4").Value = 332
$ws.Range("K4").Value = 4199
$ws.Range("L4").Value = 1424
$ws.Range("M4").Value = 2774
$ws.Range("N4").Value = 2774
$ws.Range("P4").Value = 173
$ws.Range("Q4").Value = 55
$ws.Range("R4").Value = -211
$ws.Range("S4").Value = -20
$ws.Range("T4").Value = 275
$ws.Range("U4").Value = -220
$ws.Range("V4").Value = 29
$ws.Range("W4").Value = 5.73
$ws.Range("X4").Value = 4.51
$ws.Range("Y4").Value = 12.68
$ws.Range("Z4").Value = 8.300000000000001
$ws.Range("AA4").Value = 51.35
$ws.Range("AB4").Value = 1525.26
$ws.Range("AC4").Value = 1921
$ws.Range("AD4").Value = 12.34
$ws.Range("AE4").Value = 16473
$ws.Range("AF4").Value = 1.44
$ws.Range("AG4").Value = 80
$ws.Range("AH4").Value = 0.34
$ws.Range("AI4").Value = 4.06
$ws.Range("AJ4").Value = 17261650
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 8898
$ws.Range("E5").Value = 507
$ws.Range("F5").Value = 507
$ws.Range("G5").Value = 498
$ws.Range("H5").Value = 370
$ws.Range("I5").Value = 370
$ws.Range("K5").Value = 6828
$ws.Range("L5").Value = 2691
$ws.Range("M5").Value = 4136
$ws.Range("N5").Value = 4136
$ws.Range("P5").Value = 205
$ws.Range("Q5").Value = 210
$ws.Range("R5").Value = 248
$ws.Range("S5").Value = -20
$ws.Range("T5").Value = 162
$ws.Range("U5").Value = 48
$ws.Range("V5").Value = 22
$ws.Range("W5").Value = 5.7
$ws.Range("X5").Value = 4.15
$ws.Range("Y5").Value = 10.7
$ws.Range("Z5").Value = 6.7
$ws.Range("AA5").Value = 65.06
$ws.Range("AB5").Value = 1930.19
$ws.Range("AC5").Value = 2111
$ws.Range("AD5").Value = 15.16
$ws.Range("AE5").Value = 20565
$ws.Range("AF5").Value = 1.56
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 0.31
$ws.Range("AI5").Value = 5.44
$ws.Range("AJ5").Value = 20535282
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 13517
$ws.Range("E6").Value = 481
$ws.Range("F6").Value = 481
$ws.Range("G6").Value = 504
$ws.Range("H6").Value = 389
$ws.Range("I6").Value = 389
$ws.Range("K6").Value = 6789
$ws.Range("L6").Value = 2315
$ws.Range("M6").Value = 4474
$ws.Range("N6").Value = 4474
$ws.Range("P6").Value = 205
$ws.Range("Q6").Value = 68
$ws.Range("R6").Value = -457
$ws.Range("S6").Value = -24
$ws.Range("T6").Value = 206
$ws.Range("U6").Value = -138
$ws.Range("V6").Value = 19
$ws.Range("W6").Value = 3.56
$ws.Range("X6").Value = 2.88
$ws.Range("Y6").Value = 9.029999999999999
$ws.Range("Z6").Value = 5.71
$ws.Range("AA6").Value = 51.75
$ws.Range("AB6").Value = 2094.91
$ws.Range("AC6").Value = 1893
$ws.Range("AD6").Value = 10.59
$ws.Range("AE6").Value = 22242
$ws.Range("AF6").Value = 0.9
$ws.Range("AG6").Value = 290
$ws.Range("AH6").Value = 1.45
$ws.Range("AI6").Value = 15
$ws.Range("AJ6").Value = 20535282

# Row 7
$ws.Range("D7").Value = 11996
$ws.Range("E7").Value = 279
$ws.Range("G7").Value = 300
$ws.Range("H7").Value = 229
$ws.Range("I7").Value = 229
$ws.Range("K7").Value = 7325
$ws.Range("L7").Value = 2686
$ws.Range("M7").Value = 4639
$ws.Range("N7").Value = 4649
$ws.Range("P7").Value = 207
$ws.Range("Q7").Value = 887
$ws.Range("R7").Value = -843
$ws.Range("S7").Value = -112
$ws.Range("T7").Value = 879
$ws.Range("U7").Value = 93
$ws.Range("W7").Value = 2.33
$ws.Range("X7").Value = 1.91
$ws.Range("Y7").Value = 5.02
$ws.Range("Z7").Value = 3.25
$ws.Range("AA7").Value = 57.9
$ws.Range("AC7").Value = 1115
$ws.Range("AD7").Value = 10
$ws.Range("AE7").Value = 23113
$ws.Range("AF7").Value = 0.48
$ws.Range("AG7").Value = 290
$ws.Range("AH7").Value = 2.6
$ws.Range("AI7").Value = 26

# Row 8
$ws.Range("D8").Value = 12281
$ws.Range("E8").Value = 353
$ws.Range("G8").Value = 368
$ws.Range("H8").Value = 282
$ws.Range("I8").Value = 282
$ws.Range("K8").Value = 7641
$ws.Range("L8").Value = 2778
$ws.Range("M8").Value = 4863
$ws.Range("N8").Value = 4895
$ws.Range("P8").Value = 207
$ws.Range("Q8").Value = 456
$ws.Range("R8").Value = -410
$ws.Range("S8").Value = -51
$ws.Range("T8").Value = 392
$ws.Range("U8").Value = 136
$ws.Range("W8").Value = 2.87
$ws.Range("X8").Value = 2.3
$ws.Range("Y8").Value = 5.91
$ws.Range("Z8").Value = 3.77
$ws.Range("AA8").Value = 57.11
$ws.Range("AC8").Value = 1373
$ws.Range("AD8").Value = 8.119999999999999
$ws.Range("AE8").Value = 24336
$ws.Range("AF8").Value = 0.46
$ws.Range("AG8").Value = 290
$ws.Range("AH8").Value = 2.6
$ws.Range("AI8").Value = 21.12

# Row 9
$ws.Range("D9").Value = 12616
$ws.Range("E9").Value = 413
$ws.Range("G9").Value = 425
$ws.Range("H9").Value = 325
$ws.Range("I9").Value = 325
$ws.Range("K9").Value = 7947
$ws.Range("L9").Value = 2819
$ws.Range("M9").Value = 5127
$ws.Range("N9").Value = 5196
$ws.Range("P9").Value = 207
$ws.Range("Q9").Value = 529
$ws.Range("R9").Value = -293
$ws.Range("S9").Value = -87
$ws.Range("T9").Value = 278
$ws.Range("U9").Value = 202
$ws.Range("W9").Value = 3.27
$ws.Range("X9").Value = 2.57
$ws.Range("Y9").Value = 6.44
$ws.Range("Z9").Value = 4.17
$ws.Range("AA9").Value = 54.99
$ws.Range("AC9").Value = 1581
$ws.Range("AD9").Value = 7.05
$ws.Range("AE9").Value = 25832
$ws.Range("AF9").Value = 0.43
$ws.Range("AG9").Value = 290
$ws.Range("AH9").Value = 2.6
$ws.Range("AI9").Value = 18.34
